# Mass Interview uncheck for Only Available Interviewers
# Adds the 151_fnlrgsn re-run timestamp fix plus the 152 cycle rows to the
# AMSIN quick-interview history, and the 152_betachgs row to BETA.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# AMSIN sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AMSIN")

# Row 25 already had data but the cells were missing the "general" row
# style (s=6) that the rest of the table uses, and the run-time stamp
# needs a tiny precision correction.
$ws.Range("A25:G25").Font.Name = "Calibri"
$ws.Range("B25").Value = 44476.38135983796
$ws.Range("B25").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 26 - new sprint run "152_fstcycle"
$ws.Range("A26").Value = "2021-10-26"
$ws.Range("B26").Value = 44495.65899979167
$ws.Range("B26").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C26").Value = "152_fstcycle"
$ws.Range("D26").Value = 96
$ws.Range("E26").Value = 96
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 3.04
$ws.Range("A26:G26").Font.Name = "Calibri"

# Row 27 - new sprint run "152_fnlrgrsn"
$ws.Range("A27").Value = "2021-10-28"
$ws.Range("B27").Value = 44497.39866228009
$ws.Range("B27").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C27").Value = "152_fnlrgrsn"
$ws.Range("D27").Value = 96
$ws.Range("E27").Value = 96
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 2.75
$ws.Range("A27:G27").Font.Name = "Calibri"

# ---------------------------------------------------------------------
# BETA sheet
# ---------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

# Row 15 - new sprint run "152_betachgs" (Only Available Interviewers
# unchecked => some cases fail)
$wsBeta.Range("A15").Value = "2021-10-28"
$wsBeta.Range("B15").Value = 44497.70868983068
$wsBeta.Range("B15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsBeta.Range("C15").Value = "152_betachgs"
$wsBeta.Range("D15").Value = 96
$wsBeta.Range("E15").Value = 88
$wsBeta.Range("F15").Value = 8
$wsBeta.Range("G15").Value = 6.03
